$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 2 (weekly update: newest record goes on top,
# pushing every existing row down by one; the previous last row, 34,
# becomes row 35 unchanged).
$ws.Rows(2).Insert()

# The insert copies the header row's formatting onto the new row; strip
# that back to plain (unstyled) cells like the rest of the data rows.
$ws.Range("A2:T2").ClearFormats()

# Re-apply the date/time number format used by the other rows' "Fecha" column.
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new record.
$ws.Range("A2").Value = 8
$ws.Range("B2").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C2").Value = 'Coquimbo'
$ws.Range("D2").Value = 44812
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 'Fruta'
$ws.Range("G2").Value = 100101
$ws.Range("H2").Value = 'Berries'
$ws.Range("I2").Value = 100101001
$ws.Range("J2").Value = 'Arándano (blue)'
$ws.Range("K2").Value = 'Sin especificar'
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 15500
$ws.Range("Q2").Value = '$/bandeja 2 kilos'
$ws.Range("R2").Value = 'Provincia de Limarí'
$ws.Range("S2").Value = 7750
$ws.Range("T2").Value = 2
